# Inserts a new data row at row 148 (pushing existing rows 148-188 down to
# 149-189) and populates the new row with a fresh Pomelo price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 148 and below down by one row, copying formatting (incl. the
# date-formatted style on column D) along with them.
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new record values.
$ws.Cells.Item(148, 1).Value = 10
$ws.Cells.Item(148, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(148, 3).Value = "La Araucanía"
$ws.Cells.Item(148, 4).Value = 44588
$ws.Cells.Item(148, 5).Value = 9
$ws.Cells.Item(148, 6).Value = "Fruta"
$ws.Cells.Item(148, 7).Value = 100102
$ws.Cells.Item(148, 8).Value = "Cítricos"
$ws.Cells.Item(148, 9).Value = 100102006
$ws.Cells.Item(148, 10).Value = "Pomelo"
$ws.Cells.Item(148, 11).Value = "Start Ruby"
$ws.Cells.Item(148, 12).Value = "Primera"
$ws.Cells.Item(148, 13).Value = 100
$ws.Cells.Item(148, 14).Value = 13000
$ws.Cells.Item(148, 15).Value = 13000
$ws.Cells.Item(148, 16).Value = 13000
$ws.Cells.Item(148, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(148, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(148, 19).Value = 867
$ws.Cells.Item(148, 20).Value = 15
